# "cetak pdu taruna 2021-1"
#
# This document prints two identical cadet (taruna) uniform order slips
# side by side (one Word table row, two cells). Each slip shows a set of
# mail-merge field results (name, shoe/kaos/topi size, and a row of
# measurements) that were cached as plain text in the document the last
# time the merge was run. This edit "re-prints" the slips for a new
# recipient, so every cached field result changes (some numbers repeat
# with different target values, so each one is located and replaced in
# strict left-to-right reading order rather than via a blind
# find-and-replace-all).

$d = $word.ActiveDocument

function Set-NextMatch {
    param($Range, $OldText, $NewText)

    # wdFindStop (0) so the search never wraps back past text we already
    # fixed up, and Replace = wdReplaceNone (0) so we can see exactly what
    # was matched before overwriting just that span - this keeps every
    # other occurrence of the same digits untouched.
    $found = $Range.Find.Execute(
        $OldText, $true, $false, $false, $false, $false,
        $true, 0, $false, "", 0)

    if (-not $found) {
        throw "Could not find '$OldText' while applying the reprint edit."
    }

    $Range.Text = $NewText
}

# Each slip (table cell) repeats the same sequence of cached field
# results; apply it once per slip, in document order.
$replacements = @(
    @("U1", "U10"),              # Batch/uniform code
    @("MULYONO", "ALI IQSAN S."), # Name
    @("40", "41"),                # Shoes size
    @("XL", "L"),                 # Kaos (t-shirt) size
    @("56", "58"),                # Topi (cap) size
    @("48", "46"),                # Uk. Baju measurement 1
    @("25", "24"),                # Uk. Baju measurement 2
    @("19", "18"),                # Uk. Baju measurement 3
    @("29", "28"),                # Uk. Baju measurement 4
    @("28", "26"),                # Uk. Baju measurement 5
    @("29", "27"),                # Uk. Baju measurement 6
    @("41", "40")                 # Uk. Baju measurement 8 (measurement 7, "71", is unchanged)
)

$slipCount = $d.Tables.Item(1).Rows.Count * $d.Tables.Item(1).Columns.Count
$rng = $d.Content

for ($slip = 0; $slip -lt $slipCount; $slip++) {
    foreach ($pair in $replacements) {
        Set-NextMatch $rng $pair[0] $pair[1]
    }
}
